$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.752.88'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '2.478.94'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '318.78'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '93.51'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.60%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.553'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.87%  '
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.23'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0863'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +9.24%  '
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '2.861.50'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.90'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.78'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = '2.477.69'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('E17').Value = '  +2.81%  '
$ws.Range('D18').Value = '41.715.06'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.26'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.33'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '239.99'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('E25').Value = '  +2.42%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.79'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.83'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.17'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '157.83'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.53'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.59'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0769'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.37'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.91%  '
$ws.Range('E37').Value = '  +5.22%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.93'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.05'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.47'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +10.57%  '
$ws.Range('D43').Value = '2.000.76'
$ws.Range('E43').Value = '  +2.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '19.22'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.37%  '
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.99'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.37'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.37%  '
$ws.Range('D48').Value = '2.717.73'
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '97.63'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '74.36'
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '67.25'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.17%  '
